$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 114; this shifts the existing rows 114..183
# down to 115..184 (all values/styles carried along automatically), and the
# sheet's used-range dimension grows from A1:R183 to A1:R184.
$ws.Rows.Item(114).Insert()

# Populate the newly-inserted (blank) row 114 with the new weekly record.
$ws.Range("A114").Value = 7
$ws.Range("B114").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C114").Value = 'Ñuble'
$ws.Range("D114").Value = 44582
$ws.Range("E114").Value = 16
$ws.Range("F114").Value = 100112003
$ws.Range("G114").Value = 'Ajo'
$ws.Range("H114").Value = 'Chino'
$ws.Range("I114").Value = 'Primera'
$ws.Range("J114").Value = 60
$ws.Range("K114").Value = 19000
$ws.Range("L114").Value = 20000
$ws.Range("M114").Value = 19500
$ws.Range("N114").Value = '$/caja 10 kilos'
$ws.Range("O114").Value = 'China'
$ws.Range("P114").Value = 1950
$ws.Range("Q114").Value = 10
$ws.Range("R114").Value = 'Hortaliza'
